$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (HDKSOE / 009540.KS)
$ws.Range("D2").Value = 425500
$ws.Range("E2").Value = 46
$ws.Range("F2").Value = 3.78
$ws.Range("K2").Value = 57.6
$ws.Range("N2").Value = 54.82400714602223

# Row 3 (HD HYUNDAI MIPO / 010620.KS)
$ws.Range("K3").Value = 56.4
$ws.Range("N3").Value = 54.82400714602223

# Row 4 (Hanwha Ocean / 042660.KS)
$ws.Range("D4").Value = 107100
$ws.Range("E4").Value = 18
$ws.Range("F4").Value = -0.65
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 76
$ws.Range("K4").Value = 51.4
$ws.Range("N4").Value = 54.82400714602223

# Row 5 (SamsungHvyInd / 010140.KS)
$ws.Range("K5").Value = 47.6
$ws.Range("N5").Value = 54.82400714602223
